# Move the HTTP2-logo cluster on slide 2: the whole stack of logos (plus the
# grouped "1/2/3" step-number arrows) was nudged up and slightly left.
# Every affected shape moves by the same delta: dx = -30389 EMU, dy = -511590 EMU.
#
# PowerPoint's Shape.Left/.Top are single-precision (float32) point values,
# and EMU-on-save = floor(pointValue * 12700). The literals below are the
# float32 point values (EMU/12700, nudged to the nearest representable
# float32 that still truncates back to the exact target EMU) so the saved
# XML lands exactly on the target offsets.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

$targets = @{
    "Picture 3" = @(663.13232421875, 250.2141876220703)
    "Picture 2" = @(153.54087829589844, 233.6835479736328)
    "Group 21"  = @(277.5937194824219, 263.5740966796875)
    "Picture 4" = @(498.4265441894531, 301.20538330078125)
    "Picture 6" = @(426.4071044921875, 298.7297668457031)
    "Picture 8" = @(552.7637329101562, 300.9471740722656)
}

foreach ($name in $targets.Keys) {
    $pair = $targets[$name]
    $shape = $s.Shapes.Item($name)
    $shape.Left = $pair[0]
    $shape.Top = $pair[1]
}
